$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 391
$ws.Range("F5").Value = 917
$ws.Range("F6").Value = 149
$ws.Range("F7").Value = 925
$ws.Range("F8").Value = 724
$ws.Range("F9").Value = 172
$ws.Range("F12").Value = 760
$ws.Range("F13").Value = 250
$ws.Range("F14").Value = 544
$ws.Range("F16").Value = 1278
$ws.Range("F17").Value = 112
$ws.Range("F18").Value = 418
$ws.Range("F19").Value = 1050
$ws.Range("F20").Value = 2782
$ws.Range("F21").Value = 1259
$ws.Range("F22").Value = 642
$ws.Range("F24").Value = 1235
$ws.Range("F25").Value = 52
$ws.Range("F26").Value = 962
$ws.Range("F27").Value = 315
$ws.Range("F28").Value = 533
$ws.Range("F29").Value = 1298

# --- 演出 (sheet2) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 347
$ws.Range("F7").Value = 4
$ws.Range("F10").Value = 150

# --- 本地生活 (sheet3) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 716

# --- 全部类型 (sheet4) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 716
$ws.Range("F5").Value = 391
$ws.Range("F9").Value = 347
$ws.Range("F12").Value = 917
$ws.Range("F13").Value = 149
$ws.Range("F14").Value = 4
$ws.Range("F15").Value = 925
$ws.Range("F16").Value = 724
$ws.Range("F17").Value = 172
$ws.Range("F21").Value = 150
$ws.Range("F25").Value = 760
$ws.Range("F26").Value = 250
$ws.Range("F27").Value = 544
$ws.Range("F29").Value = 1278
$ws.Range("F30").Value = 112
$ws.Range("F31").Value = 418
$ws.Range("F32").Value = 1050
$ws.Range("F33").Value = 2782
$ws.Range("F34").Value = 1259
$ws.Range("F35").Value = 642
$ws.Range("F37").Value = 1235
$ws.Range("F38").Value = 52
$ws.Range("F40").Value = 962
$ws.Range("F41").Value = 315
$ws.Range("F42").Value = 534
$ws.Range("F43").Value = 1298
